$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new weekly data row at row 38, pushing the existing
# rows 38-108 down to 39-109 (dimension grows from A1:R108 to A1:R109).
$ws.Rows("38").Insert()

$ws.Range("A38").Value = 4
$ws.Range("B38").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C38").Value = "Los Lagos"
$ws.Range("D38").Value = 44469
$ws.Range("E38").Value = 10
$ws.Range("F38").Value = 100112009
$ws.Range("G38").Value = "Acelga"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 100
$ws.Range("K38").Value = 4000
$ws.Range("L38").Value = 4000
$ws.Range("M38").Value = 4000
$ws.Range("N38").Value = "$/docena de atados (4 kilos)"
$ws.Range("O38").Value = "Región del Maule"
$ws.Range("P38").Value = 1000
$ws.Range("Q38").Value = 4
$ws.Range("R38").Value = "Hortaliza"
